$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25: new "Suite du dossier de projet" task entry, 2 periods
$ws.Range("C25").Value = "Suite du dossier de projet"
$ws.Range("D25").Value = 2
$ws.Range("D25").HorizontalAlignment = -4108
$ws.Range("D25").VerticalAlignment = -4108

# Row 26: reuse existing "Recherches sur les réseaux neuronaux" task text, 2 periods
$ws.Range("C26").Value = "Recherches sur les réseaux neuronaux"
$ws.Range("D26").Value = 2
$ws.Range("D26").HorizontalAlignment = -4108
$ws.Range("D26").VerticalAlignment = -4108

# Row 27: new "Mise à jour des Uses cases" task entry, 2 periods
$ws.Range("C27").Value = "Mise à jour des Uses cases"
$ws.Range("D27").Value = 2
$ws.Range("D27").HorizontalAlignment = -4108
$ws.Range("D27").VerticalAlignment = -4108

# Move the active selection to C28 (as recorded in the saved view state)
[void]$ws.Range("C28").Select()
